# feat: add 2022-Q4 data
#
# 1. "总计" summary sheet: insert the new 2022-Q4 totals as the new row 2,
#    pushing the former row 2 (2022-Q3 totals) down to row 3 (its index
#    bumps from 0 to 1).
# 2. The former "2022-Q3" detail sheet keeps its data but is copied into a
#    brand-new sheet placed right after it (still named "2022-Q3"), and the
#    original sheet slot is renamed "2022-Q4" and repopulated with the new
#    quarter's fund holdings.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: "总计" sheet - shift old Q3 summary row down, add new Q4 row
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

# Write the (previously row 2) Q3 summary into row 3, with its index bumped.
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 1
$wsTotal.Range("D3").Value = 0.01
# Match the formatting of the index column (A2 already carries it).
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)

# Overwrite row 2 with the new Q4 summary values.
$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 8
$wsTotal.Range("D2").Value = 0.54

# ---------------------------------------------------------------------
# Step 2: duplicate the old "2022-Q3" detail sheet's content into a new
# sheet (keeping the "2022-Q3" name), then repurpose the original sheet
# slot as "2022-Q4" and repopulate it with the new fund table.
# ---------------------------------------------------------------------
$wsQ3Old = $wb.Worksheets.Item(2)

$wsQ3New = $wb.Worksheets.Add($null, $wsQ3Old)
$wsQ3Old.Range("A1:H2").Copy($wsQ3New.Range("A1"))
$wsQ3New.Range("A1").ClearContents()

# Free up the "2022-Q3" name on the old sheet before claiming it on the new one.
$wsQ3Old.Name = "2022-Q4"
$wsQ3New.Name = "2022-Q3"

# Clear the old single-fund row so we can lay down the full Q4 table.
$wsQ3Old.Range("A1:H2").ClearContents()

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = [char]([int][char]'B' + $i)
    $wsQ3Old.Range("$col" + "1").Value = $headers[$i]
}
# Copy header formatting (s=2 style) from the 总计 sheet's header row.
$wsTotal.Range("B1").Copy()
$wsQ3Old.Range("B1:H1").PasteSpecial(-4122)

$rows = @(
    @("160921", "大成多策略混合（LOF）A", "8.25", "87.54", "3.65", "0.3011", 8),
    @("016062", "大成多策略混合（LOF）C", "6.19", "87.54", "3.65", "0.2259", 8),
    @("013166", "东兴宸祥量化混合A", "0.38", "93.88", "1.36", "0.0052", 2),
    @("009327", "东兴兴晟混合A", "0.38", "79.79", "1.15", "0.0044", 3),
    @("013167", "东兴宸祥量化混合C", "0.08", "93.88", "1.36", "0.0011", 2),
    @("000926", "中信建投睿信灵活配置混合A", "0.10", "83.25", "1.03", "0.0010", 8),
    @("009328", "东兴兴晟混合C", "0.08", "79.79", "1.15", "0.0009", 3),
    @("004676", "中信建投睿信灵活配置混合C", "0.03", "83.25", "1.03", "0.0003", 8)
)

# Stage the numeric-looking text columns (fund code, scale, position,
# ratio, market value) in a scratch area formatted as Text so Excel
# doesn't coerce them (and strip leading zeros) into real numbers, then
# paste-values them into place, which drops the Text number format
# again so the final cells carry no explicit style - matching the
# inlineStr cells produced by the original export.
$scratchCode = $wsTotal.Range("Z1:Z8")
$scratchCode.NumberFormat = "@"
$scratchScale = $wsTotal.Range("AA1:AA8")
$scratchScale.NumberFormat = "@"
$scratchPos = $wsTotal.Range("AB1:AB8")
$scratchPos.NumberFormat = "@"
$scratchRatio = $wsTotal.Range("AC1:AC8")
$scratchRatio.NumberFormat = "@"
$scratchValue = $wsTotal.Range("AD1:AD8")
$scratchValue.NumberFormat = "@"

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 1
    $data = $rows[$i]
    $wsTotal.Range("Z$r").Value = $data[0]
    $wsTotal.Range("AA$r").Value = $data[2]
    $wsTotal.Range("AB$r").Value = $data[3]
    $wsTotal.Range("AC$r").Value = $data[4]
    $wsTotal.Range("AD$r").Value = $data[5]
}

$scratchCode.Copy()
$wsQ3Old.Range("B2:B9").PasteSpecial(-4163)
$scratchScale.Copy()
$wsQ3Old.Range("D2:D9").PasteSpecial(-4163)
$scratchPos.Copy()
$wsQ3Old.Range("E2:E9").PasteSpecial(-4163)
$scratchRatio.Copy()
$wsQ3Old.Range("F2:F9").PasteSpecial(-4163)
$scratchValue.Copy()
$wsQ3Old.Range("G2:G9").PasteSpecial(-4163)

# Clean up the scratch area.
$wsTotal.Range("Z1:AD8").Clear()

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $data = $rows[$i]
    $wsQ3Old.Range("A$r").Value = $i
    $wsQ3Old.Range("C$r").Value = $data[1]
    $wsQ3Old.Range("H$r").Value = $data[6]
}
# Copy the index-column formatting (s=2 style) down too.
$wsTotal.Range("A2").Copy()
$wsQ3Old.Range("A2:A9").PasteSpecial(-4122)

Write-Host "done"
